# Update "想去人数" (attendance interest count) figures in the F column
# across the "展览" and "全部类型" worksheets, matching the freshly
# regenerated gh-pages data output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F
$updates = @{
    2  = 2968
    4  = 104
    5  = 6726
    6  = 1712
    10 = 118
    11 = 25
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
